$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: set values
$ws.Range("D2").Value = 16333
$ws.Range("E2").Value = -1037
$ws.Range("F2").Value = -1037
$ws.Range("G2").Value = -1406
$ws.Range("H2").Value = -1111
$ws.Range("I2").Value = -1115
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 12383
$ws.Range("L2").Value = 9366
$ws.Range("M2").Value = 3017
$ws.Range("N2").Value = 2568
$ws.Range("O2").Value = 449
$ws.Range("P2").Value = 447
$ws.Range("Q2").Value = -9
$ws.Range("R2").Value = 412
$ws.Range("S2").Value = -953
$ws.Range("T2").Value = 74
$ws.Range("U2").Value = -83
$ws.Range("V2").Value = 4379
$ws.Range("W2").Value = -6.35
$ws.Range("X2").Value = -6.8
$ws.Range("Y2").Value = -35.34
$ws.Range("Z2").Value = -8.17
$ws.Range("AA2").Value = 310.44
$ws.Range("AB2").Value = 469.11
$ws.Range("AC2").Value = -12490
$ws.Range("AD2").Value = -0.76
$ws.Range("AE2").Value = 28990
$ws.Range("AF2").Value = 0.33
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 8930907

# Row 3: set values
$ws.Range("D3").Value = 15222
$ws.Range("E3").Value = 361
$ws.Range("F3").Value = 361
$ws.Range("G3").Value = 212
$ws.Range("H3").Value = 169
$ws.Range("I3").Value = 157
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 12247
$ws.Range("L3").Value = 9068
$ws.Range("M3").Value = 3179
$ws.Range("N3").Value = 2750
$ws.Range("O3").Value = 429
$ws.Range("P3").Value = 447
$ws.Range("Q3").Value = 1396
$ws.Range("R3").Value = -237
$ws.Range("S3").Value = -639
$ws.Range("T3").Value = 53
$ws.Range("U3").Value = 1342
$ws.Range("V3").Value = 3929
$ws.Range("W3").Value = 2.37
$ws.Range("X3").Value = 1.11
$ws.Range("Y3").Value = 5.89
$ws.Range("Z3").Value = 1.37
$ws.Range("AA3").Value = 285.22
$ws.Range("AB3").Value = 509.95
$ws.Range("AC3").Value = 1754
$ws.Range("AD3").Value = 5.67
$ws.Range("AE3").Value = 31051
$ws.Range("AF3").Value = 0.32
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 8930907

# Row 4: set values
$ws.Range("D4").Value = 19585
$ws.Range("E4").Value = 776
$ws.Range("F4").Value = 776
$ws.Range("G4").Value = 362
$ws.Range("H4").Value = 186
$ws.Range("I4").Value = 166
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 15134
$ws.Range("L4").Value = 11840
$ws.Range("M4").Value = 3294
$ws.Range("N4").Value = 2885
$ws.Range("O4").Value = 410
$ws.Range("P4").Value = 447
$ws.Range("Q4").Value = 1260
$ws.Range("R4").Value = -358
$ws.Range("S4").Value = -1
$ws.Range("T4").Value = 73
$ws.Range("U4").Value = 1187
$ws.Range("V4").Value = 4165
$ws.Range("W4").Value = 3.96
$ws.Range("X4").Value = 0.95
$ws.Range("Y4").Value = 5.91
$ws.Range("Z4").Value = 1.36
$ws.Range("AA4").Value = 359.44
$ws.Range("AB4").Value = 538.9299999999999
$ws.Range("AC4").Value = 1864
$ws.Range("AD4").Value = 7.56
$ws.Range("AE4").Value = 32565
$ws.Range("AF4").Value = 0.43
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 8930907
# Row 4: clear cells
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5: set values
$ws.Range("D5").Value = 22408
$ws.Range("E5").Value = 1150
$ws.Range("F5").Value = 1150
$ws.Range("G5").Value = 796
$ws.Range("H5").Value = 565
$ws.Range("I5").Value = 523
$ws.Range("J5").Value = 42
$ws.Range("K5").Value = 15829
$ws.Range("L5").Value = 12006
$ws.Range("M5").Value = 3823
$ws.Range("N5").Value = 3380
$ws.Range("O5").Value = 444
$ws.Range("P5").Value = 447
$ws.Range("Q5").Value = -584
$ws.Range("R5").Value = -518
$ws.Range("S5").Value = 414
$ws.Range("T5").Value = 104
$ws.Range("U5").Value = -688
$ws.Range("V5").Value = 5012
$ws.Range("W5").Value = 5.13
$ws.Range("X5").Value = 2.52
$ws.Range("Y5").Value = 16.7
$ws.Range("Z5").Value = 3.65
$ws.Range("AA5").Value = 314.02
$ws.Range("AB5").Value = 648.4
$ws.Range("AC5").Value = 5856
$ws.Range("AD5").Value = 2.95
$ws.Range("AE5").Value = 38153
$ws.Range("AF5").Value = 0.45
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 8930907
# Row 5: clear cells
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6: set values
$ws.Range("D6").Value = 22927
$ws.Range("E6").Value = 1537
$ws.Range("F6").Value = 1537
$ws.Range("G6").Value = 1262
$ws.Range("H6").Value = 945
$ws.Range("I6").Value = 892
$ws.Range("K6").Value = 17262
$ws.Range("L6").Value = 12864
$ws.Range("M6").Value = 4398
$ws.Range("N6").Value = 3914
$ws.Range("P6").Value = 447
$ws.Range("Q6").Value = 1041
$ws.Range("R6").Value = -98
$ws.Range("S6").Value = 13
$ws.Range("T6").Value = 85
$ws.Range("U6").Value = 956
$ws.Range("V6").Value = 5301
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 4.12
$ws.Range("Y6").Value = 24.46
$ws.Range("Z6").Value = 5.71
$ws.Range("AA6").Value = 292.48
$ws.Range("AB6").Value = 778.27
$ws.Range("AC6").Value = 9987
$ws.Range("AD6").Value = 2.18
$ws.Range("AE6").Value = 44186
$ws.Range("AF6").Value = 0.49
$ws.Range("AI6").Value = 4.97
$ws.Range("AJ6").Value = 8930907
# Row 6: clear cells
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7: set values
$ws.Range("D7").Value = 21369
$ws.Range("E7").Value = 1188
$ws.Range("G7").Value = 898
$ws.Range("I7").Value = 640
$ws.Range("W7").Value = 5.56
$ws.Range("AC7").Value = 7166
$ws.Range("AD7").Value = 2.49
# Row 7: clear cells
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: set values
$ws.Range("D8").Value = 21276
$ws.Range("E8").Value = 960
$ws.Range("G8").Value = 720
$ws.Range("I8").Value = 473
$ws.Range("W8").Value = 4.51
$ws.Range("AC8").Value = 5296
$ws.Range("AD8").Value = 3.37
# Row 8: clear cells
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: clear cells
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
